$wb = $excel.ActiveWorkbook

# --- Sheet "TestCasesFlag": change selection to whole used range A1:B6 (no explicit active cell) ---
$wsFlag = $wb.Worksheets.Item("TestCasesFlag")
$wsFlag.Activate()
$wsFlag.Range("A1:B6").Select() | Out-Null

# --- Sheet "Data": rotate the 3 data rows (row2 -> row4, row3 -> row2, row4 -> row3) ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate()

# Capture current values of the three data rows (columns A..H)
$row2 = $wsData.Range("A2:H2").Value2
$row3 = $wsData.Range("A3:H3").Value2
$row4 = $wsData.Range("A4:H4").Value2

# Write back in rotated order: old row3 -> row2, old row4 -> row3, old row2 -> row4
$wsData.Range("A2:H2").Value2 = $row3
$wsData.Range("A3:H3").Value2 = $row4
$wsData.Range("A4:H4").Value2 = $row2

# Update the selection/active cell on the Data sheet
$wsData.Range("H9").Select() | Out-Null
